$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.158.81"
$ws.Range("E2").Value = "  -1.92%  "
$ws.Range("D3").Value = "1.822.80"
$ws.Range("E3").Value = "  -1.35%  "
$ws.Range("E4").Value = "  -0.52%  "
$ws.Range("D5").Value = "312.87"
$ws.Range("E5").Value = "  -1.97%  "
$ws.Range("E6").Value = "  -0.56%  "
$ws.Range("D7").Value = "0.4231"
$ws.Range("E7").Value = "  -1.86%  "
$ws.Range("E8").Value = "  -1.57%  "
$ws.Range("D9").Value = "0.07246"
$ws.Range("E9").Value = "  -1.47%  "
$ws.Range("D10").Value = "0.8580"
$ws.Range("E10").Value = "  -2.61%  "
$ws.Range("D11").Value = "20.96"
$ws.Range("E11").Value = "  -2.90%  "
$ws.Range("D12").Value = "1.837.85"
$ws.Range("E12").Value = "  -0.73%  "
$ws.Range("D13").Value = "6.710"
$ws.Range("E13").Value = "  -0.42%  "
$ws.Range("D14").Value = "0.07089"
$ws.Range("E14").Value = "  -0.22%  "
$ws.Range("D15").Value = "5.304"
$ws.Range("E15").Value = "  -2.85%  "
$ws.Range("D16").Value = "89.58"
$ws.Range("E16").Value = "  +2.10%  "
$ws.Range("E17").Value = "  -0.73%  "
$ws.Range("D18").Value = "0.000008850"
$ws.Range("E18").Value = "  -1.63%  "
$ws.Range("E19").Value = "  -0.63%  "
$ws.Range("E20").Value = "  -3.01%  "
$ws.Range("D21").Value = "27.230.94"
$ws.Range("E21").Value = "  -1.68%  "
$ws.Range("D22").Value = "5.130"
$ws.Range("E22").Value = "  -2.43%  "
$ws.Range("E23").Value = "  -2.46%  "
$ws.Range("D24").Value = "2.052.24"
$ws.Range("E24").Value = "  -1.28%  "
$ws.Range("E25").Value = "  -2.34%  "
$ws.Range("D26").Value = "152.48"
$ws.Range("E26").Value = "  -1.99%  "
$ws.Range("D27").Value = "2.188"
$ws.Range("E27").Value = "  +2.24%  "
$ws.Range("D28").Value = "18.40"
$ws.Range("E28").Value = "  -1.09%  "
$ws.Range("D29").Value = "5.224"
$ws.Range("E29").Value = "  -3.06%  "
$ws.Range("D30").Value = "116.40"
$ws.Range("E30").Value = "  -3.11%  "
$ws.Range("D31").Value = "0.08843"
$ws.Range("E31").Value = "  -0.86%  "
$ws.Range("D32").Value = "1.189"
$ws.Range("D33").Value = "0.7498"
$ws.Range("E33").Value = "  -3.57%  "
$ws.Range("D34").Value = "4.438"
$ws.Range("E34").Value = "  -2.72%  "
$ws.Range("D35").Value = "2.832"
$ws.Range("E35").Value = "  -2.73%  "
$ws.Range("D36").Value = "1.006"
$ws.Range("E36").Value = "  -0.61%  "
$ws.Range("D37").Value = "1.117"
$ws.Range("E37").Value = "  -2.08%  "
$ws.Range("D38").Value = "0.01970"
$ws.Range("E38").Value = "  -0.09%  "
$ws.Range("E39").Value = "  -1.73%  "
$ws.Range("D40").Value = "7.298"
$ws.Range("E40").Value = "  +1.06%  "
$ws.Range("D41").Value = "2.879"
$ws.Range("E41").Value = "  +0.54%  "
$ws.Range("D42").Value = "0.1692"
$ws.Range("E42").Value = "  +0.76%  "
$ws.Range("E43").Value = "  -2.65%  "
$ws.Range("D44").Value = "8.662"
$ws.Range("E44").Value = "  -3.30%  "
$ws.Range("D45").Value = "10.62"
$ws.Range("E45").Value = "  -0.37%  "
$ws.Range("D46").Value = "106.53"
$ws.Range("E46").Value = "  -3.79%  "
$ws.Range("D47").Value = "0.4739"
$ws.Range("E47").Value = "  +0.09%  "
$ws.Range("D48").Value = "1.006"
$ws.Range("E48").Value = "  -0.62%  "
$ws.Range("D49").Value = "0.06390"
$ws.Range("E49").Value = "  -1.68%  "
$ws.Range("E50").Value = "  -2.10%  "
$ws.Range("D51").Value = "1.854"
$ws.Range("E51").Value = "  -1.95%  "
